# Add files via upload
# The "Image" column (column A) on the active sheet holds file names such as
# "shape_0.png". The images were moved into a "shape" subfolder, so every
# value in that column needs to be updated to "shape/shape_0.png" (etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^shape_\d+\.png$') {
        $cell.Value2 = "shape/" + $val
    }
}
